$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 355.625
$ws.Range("I12").Value = 363.57144
$ws.Range("K12").Value = 363.57144
$ws.Range("M12").Value = -193.57144
$ws.Range("H55").Value = 388.73334
$ws.Range("I55").Value = 196
$ws.Range("J55").Value = 500.3158
$ws.Range("K55").Value = 196
$ws.Range("L55").Value = 500.3158
$ws.Range("M55").Value = 18
$ws.Range("N55").Value = -928.3158000000001
$ws.Range("H80").Value = 378.64285
$ws.Range("I80").Value = 187.625
$ws.Range("J80").Value = 633.3333
$ws.Range("K80").Value = 562.875
$ws.Range("L80").Value = 1899.9999
$ws.Range("M80").Value = 435.125
$ws.Range("N80").Value = -3895.9999
$ws.Range("H83").Value = 378.64285
$ws.Range("I83").Value = 187.625
$ws.Range("J83").Value = 633.3333
$ws.Range("K83").Value = 1688.625
$ws.Range("L83").Value = 5699.9997
$ws.Range("M83").Value = 3303.375
$ws.Range("N83").Value = -15683.9997
$ws.Range("H88").Value = 618309.6
$ws.Range("I88").Value = 1389371.4
$ws.Range("J88").Value = 1460.2
$ws.Range("K88").Value = 1389371.4
$ws.Range("L88").Value = 1460.2
$ws.Range("M88").Value = -1388965.4
$ws.Range("N88").Value = -2272.2
$ws.Range("H91").Value = 618309.6
$ws.Range("I91").Value = 1389371.4
$ws.Range("J91").Value = 1460.2
$ws.Range("K91").Value = 1389371.4
$ws.Range("L91").Value = 1460.2
$ws.Range("M91").Value = -1387967.4
$ws.Range("N91").Value = -4268.2
$ws.Range("H132").Value = 2547.15
$ws.Range("I132").Value = 2254.0908
$ws.Range("J132").Value = 2905.3333
$ws.Range("K132").Value = 6762.2724
$ws.Range("L132").Value = 8715.999899999999
$ws.Range("M132").Value = -4232.2724
$ws.Range("N132").Value = -13775.9999

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 50000
$ws.Range("J13").Value = 50000
$ws.Range("L13").Value = 50000
$ws.Range("N13").Value = -50336
$ws.Range("H82").Value = 11627
$ws.Range("I82").Value = 4833
$ws.Range("J82").Value = 25215
$ws.Range("K82").Value = 4833
$ws.Range("L82").Value = 25215
$ws.Range("M82").Value = -4450
$ws.Range("N82").Value = -25981
$ws.Range("H85").Value = 11627
$ws.Range("I85").Value = 4833
$ws.Range("J85").Value = 25215
$ws.Range("K85").Value = 4833
$ws.Range("L85").Value = 25215
$ws.Range("M85").Value = -3507
$ws.Range("N85").Value = -27867

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 45.866665
$ws.Range("I7").Value = 38.636364
$ws.Range("J7").Value = 65.75
$ws.Range("K7").Value = 38.636364
$ws.Range("L7").Value = 65.75
$ws.Range("M7").Value = 74.363636
$ws.Range("N7").Value = -291.75
$ws.Range("H80").Value = 15875
$ws.Range("J80").Value = 15875
$ws.Range("L80").Value = 15875
$ws.Range("N80").Value = -18121
$ws.Range("H83").Value = 15875
$ws.Range("J83").Value = 15875
$ws.Range("L83").Value = 47625
$ws.Range("N83").Value = -58857
$ws.Range("H115").Value = 30500
$ws.Range("J115").Value = 30500
$ws.Range("L115").Value = 30500
$ws.Range("N115").Value = -32850
$ws.Range("H132").Value = 3288.6667
$ws.Range("I132").Value = 2266.6667
$ws.Range("J132").Value = 5332.6665
$ws.Range("K132").Value = 6800.000100000001
$ws.Range("L132").Value = 15997.9995
$ws.Range("M132").Value = -4270.000100000001
$ws.Range("N132").Value = -21057.9995
$ws.Range("H134").Value = 1890.409
$ws.Range("I134").Value = 1497.8572
$ws.Range("J134").Value = 2577.375
$ws.Range("K134").Value = 4493.571599999999
$ws.Range("L134").Value = 7732.125
$ws.Range("M134").Value = -1958.571599999999
$ws.Range("N134").Value = -12802.125

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1525.0625
$ws.Range("I4").Value = 242.78572
$ws.Range("K4").Value = 728.35716
$ws.Range("M4").Value = -616.35716
$ws.Range("H5").Value = 480.35715
$ws.Range("I5").Value = 460
$ws.Range("J5").Value = 602.5
$ws.Range("K5").Value = 1380
$ws.Range("L5").Value = 1807.5
$ws.Range("M5").Value = -1268
$ws.Range("N5").Value = -2031.5
$ws.Range("H135").Value = 480.35715
$ws.Range("I135").Value = 460
$ws.Range("J135").Value = 602.5
$ws.Range("K135").Value = 4140
$ws.Range("L135").Value = 5422.5
$ws.Range("M135").Value = -1605
$ws.Range("N135").Value = -10492.5

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 22000
$ws.Range("J15").Value = 22000
$ws.Range("L15").Value = 22000
$ws.Range("N15").Value = -22576
$ws.Range("H81").Value = 22000
$ws.Range("J81").Value = 22000
$ws.Range("L81").Value = 22000
$ws.Range("N81").Value = -23996
$ws.Range("H84").Value = 22000
$ws.Range("J84").Value = 22000
$ws.Range("L84").Value = 66000
$ws.Range("N84").Value = -75984
$ws.Range("H102").Value = 2769.875
$ws.Range("I102").Value = 2405.6667
$ws.Range("J102").Value = 3862.5
$ws.Range("K102").Value = 2405.6667
$ws.Range("L102").Value = 3862.5
$ws.Range("M102").Value = -783.6667000000002
$ws.Range("N102").Value = -7106.5
$ws.Range("H103").Value = 19000.5
$ws.Range("J103").Value = 19000.5
$ws.Range("L103").Value = 19000.5
$ws.Range("N103").Value = -21344.5

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 13425
$ws.Range("J42").Value = 13425
$ws.Range("L42").Value = 13425
$ws.Range("N42").Value = -14551
$ws.Range("H49").Value = 13425
$ws.Range("J49").Value = 13425
$ws.Range("L49").Value = 13425
$ws.Range("N49").Value = -13719
$ws.Range("H117").Value = 50000
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 50000
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 50000
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -59178

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 29294.445
$ws.Range("J75").Value = 29294.445
$ws.Range("L75").Value = 29294.445
$ws.Range("N75").Value = -31166.445
$ws.Range("H78").Value = 29294.445
$ws.Range("J78").Value = 29294.445
$ws.Range("L78").Value = 87883.33499999999
$ws.Range("N78").Value = -97243.33499999999
$ws.Range("H86").Value = 10325
$ws.Range("J86").Value = 10325
$ws.Range("L86").Value = 10325
$ws.Range("N86").Value = -12571
$ws.Range("H89").Value = 10325
$ws.Range("J89").Value = 10325
$ws.Range("L89").Value = 51625
$ws.Range("N89").Value = -62857
$ws.Range("H118").Value = 48000
$ws.Range("J118").Value = 48000
$ws.Range("L118").Value = 48000
$ws.Range("N118").Value = -51314
$ws.Range("H126").Value = 916.35297
$ws.Range("I126").Value = 617.3
$ws.Range("J126").Value = 1343.5714
$ws.Range("K126").Value = 1851.9
$ws.Range("L126").Value = 4030.7142
$ws.Range("M126").Value = 618.1000000000001
$ws.Range("N126").Value = -8970.7142
